$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (ID) holds text values throughout the sheet; force text formatting
# on the range before writing so the IDs are not reinterpreted as numbers, then
# clear the temporary number-format override so no extra styling is introduced.
$idRange = $ws.Range("A492:A532")
$idRange.NumberFormat = "@"

$ws.Range("A492").Value = "34376"
$ws.Range("J492").Value = "No"
$ws.Range("K492").Value = "Public"
$ws.Range("A493").Value = "70822"
$ws.Range("J493").Value = "No"
$ws.Range("K493").Value = "Public"
$ws.Range("A494").Value = "36236"
$ws.Range("J494").Value = "No"
$ws.Range("K494").Value = "Private"
$ws.Range("A495").Value = "25226"
$ws.Range("J495").Value = "Yes"
$ws.Range("K495").Value = "Self-employed"
$ws.Range("A496").Value = "2182"
$ws.Range("J496").Value = "No"
$ws.Range("K496").Value = "Public"
$ws.Range("A497").Value = "67177"
$ws.Range("J497").Value = "No"
$ws.Range("K497").Value = "Private"
$ws.Range("A498").Value = "39373"
$ws.Range("J498").Value = "Yes"
$ws.Range("K498").Value = "Self-employed"
$ws.Range("A499").Value = "11974"
$ws.Range("J499").Value = "No"
$ws.Range("K499").Value = "Public"
$ws.Range("A500").Value = "28645"
$ws.Range("J500").Value = "Yes"
$ws.Range("K500").Value = "Private"
$ws.Range("A501").Value = "56681"
$ws.Range("J501").Value = "Yes"
$ws.Range("K501").Value = "Private"
$ws.Range("A502").Value = "56546"
$ws.Range("J502").Value = "Yes"
$ws.Range("K502").Value = "Private"
$ws.Range("A503").Value = "21408"
$ws.Range("J503").Value = "Yes"
$ws.Range("K503").Value = "Self-employed"
$ws.Range("A504").Value = "66400"
$ws.Range("J504").Value = "Yes"
$ws.Range("K504").Value = "Retired"
$ws.Range("A505").Value = "71038"
$ws.Range("J505").Value = "No"
$ws.Range("K505").Value = "Public"
$ws.Range("A506").Value = "46785"
$ws.Range("J506").Value = "Yes"
$ws.Range("K506").Value = "Private"
$ws.Range("A507").Value = "50931"
$ws.Range("J507").Value = "No"
$ws.Range("K507").Value = "Public"
$ws.Range("A508").Value = "58978"
$ws.Range("J508").Value = "Yes"
$ws.Range("K508").Value = "Private"
$ws.Range("A509").Value = "11091"
$ws.Range("J509").Value = "Yes"
$ws.Range("K509").Value = "Private"
$ws.Range("A510").Value = "32503"
$ws.Range("J510").Value = "No"
$ws.Range("K510").Value = "Self-employed"
$ws.Range("A511").Value = "4651"
$ws.Range("J511").Value = "Yes"
$ws.Range("K511").Value = "Private"
$ws.Range("A512").Value = "54385"
$ws.Range("J512").Value = "No"
$ws.Range("K512").Value = "Self-employed"
$ws.Range("A513").Value = "25774"
$ws.Range("J513").Value = "Yes"
$ws.Range("K513").Value = "Private"
$ws.Range("A514").Value = "17718"
$ws.Range("J514").Value = "Yes"
$ws.Range("K514").Value = "Retired"
$ws.Range("A515").Value = "43054"
$ws.Range("J515").Value = "Yes"
$ws.Range("K515").Value = "Private"
$ws.Range("A516").Value = "42072"
$ws.Range("J516").Value = "Yes"
$ws.Range("K516").Value = "Self-employed"
$ws.Range("A517").Value = "64908"
$ws.Range("J517").Value = "Yes"
$ws.Range("K517").Value = "Private"
$ws.Range("A518").Value = "67432"
$ws.Range("J518").Value = "Yes"
$ws.Range("K518").Value = "Private"
$ws.Range("A519").Value = "18587"
$ws.Range("J519").Value = "Yes"
$ws.Range("K519").Value = "Self-employed"
$ws.Range("A520").Value = "38165"
$ws.Range("J520").Value = "No"
$ws.Range("K520").Value = "Private"
$ws.Range("A521").Value = "54375"
$ws.Range("J521").Value = "Yes"
$ws.Range("K521").Value = "Public"
$ws.Range("A522").Value = "58631"
$ws.Range("J522").Value = "Yes"
$ws.Range("K522").Value = "Self-employed"
$ws.Range("A523").Value = "26325"
$ws.Range("J523").Value = "Yes"
$ws.Range("K523").Value = "Private"
$ws.Range("A524").Value = "27832"
$ws.Range("J524").Value = "Yes"
$ws.Range("K524").Value = "Private"
$ws.Range("A525").Value = "6118"
$ws.Range("J525").Value = "No"
$ws.Range("K525").Value = "Private"
$ws.Range("A526").Value = "69551"
$ws.Range("J526").Value = "Yes"
$ws.Range("K526").Value = "Public"
$ws.Range("A527").Value = "16371"
$ws.Range("J527").Value = "Yes"
$ws.Range("K527").Value = "Self-employed"
$ws.Range("A528").Value = "31179"
$ws.Range("J528").Value = "Yes"
$ws.Range("K528").Value = "Private"
$ws.Range("A529").Value = "30456"
$ws.Range("J529").Value = "No"
$ws.Range("K529").Value = "Private"
$ws.Range("A530").Value = "68627"
$ws.Range("J530").Value = "Yes"
$ws.Range("K530").Value = "Private"
$ws.Range("A531").Value = "30468"
$ws.Range("J531").Value = "No"
$ws.Range("K531").Value = "Public"
$ws.Range("A532").Value = "48796"
$ws.Range("J532").Value = "Yes"
$ws.Range("K532").Value = "Private"

$idRange.ClearFormats()
